$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.456.42'
$ws.Range('E2').Value = '  +0.35%  '
$ws.Range('D3').Value = '1.942.82'
$ws.Range('E3').Value = '  -1.13%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = "'243.70"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.11%  '
$ws.Range('D6').Value = "'0.614"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.42%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = "'57.56"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.56%  '
$ws.Range('D9').Value = "'0.361"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.28%  '
$ws.Range('D10').Value = "'0.0848"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.89%  '
$ws.Range('E11').Value = '  -0.76%  '
$ws.Range('D12').Value = '2.226.31'
$ws.Range('E12').Value = '  -1.00%  '
$ws.Range('B13').Value = 'Avalanche'
$ws.Range('C13').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D13').Value = "'21.40"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -4.09%  '
$ws.Range('B14').Value = 'Polygon'
$ws.Range('C14').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D14').Value = "'0.813"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.63%  '
$ws.Range('D15').Value = "'13.48"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.10%  '
$ws.Range('E16').Value = '  -3.59%  '
$ws.Range('D17').Value = '1.952.19'
$ws.Range('E17').Value = '  -0.45%  '
$ws.Range('D18').Value = '36.412.46'
$ws.Range('E18').Value = '  +0.50%  '
$ws.Range('D19').Value = "'69.35"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.53%  '
$ws.Range('D20').Value = '0.0₃0865'
$ws.Range('E20').Value = '  -2.61%  '
$ws.Range('D21').Value = "'228.85"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.07%  '
$ws.Range('D22').Value = "'5.00"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.80%  '
$ws.Range('E23').Value = '  -0.20%  '
$ws.Range('D24').Value = "'2.37"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -5.87%  '
$ws.Range('D25').Value = "'2.30"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.34%  '
$ws.Range('D26').Value = "'9.21"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -4.14%  '
$ws.Range('D27').Value = "'161.85"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.77%  '
$ws.Range('D28').Value = "'0.135"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +6.76%  '
$ws.Range('D29').Value = "'19.21"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.74%  '
$ws.Range('E30').Value = '  -0.74%  '
$ws.Range('E31').Value = '  -4.30%  '
$ws.Range('D32').Value = "'4.59"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.44%  '
$ws.Range('D33').Value = "'0.0618"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.87%  '
$ws.Range('B34').Value = 'InternetComputer(DFINITY)'
$ws.Range('C34').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D34').Value = "'4.19"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.97%  '
$ws.Range('B35').Value = 'THORChain'
$ws.Range('C35').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D35').Value = "'6.23"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +4.74%  '
$ws.Range('E36').Value = '  -0.05%  '
$ws.Range('E37').Value = '  -1.22%  '
$ws.Range('D38').Value = "'2.18"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.47%  '
$ws.Range('D39').Value = "'3.18"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +8.78%  '
$ws.Range('D40').Value = "'0.0985"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.65%  '
$ws.Range('E41').Value = '  +0.44%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').Value = "'1.15"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.84%  '
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').Value = "'0.0209"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.64%  '
$ws.Range('D44').Value = "'16.00"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.03%  '
$ws.Range('D45').Value = '1.343.87'
$ws.Range('E45').Value = '  -0.31%  '
$ws.Range('E46').Value = '  -3.02%  '
$ws.Range('D47').Value = "'7.21"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.80%  '
$ws.Range('D48').Value = "'86.43"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.68%  '
$ws.Range('E49').Value = '  -0.05%  '
$ws.Range('D50').Value = '2.118.18'
$ws.Range('E50').Value = '  -0.90%  '
$ws.Range('D51').Value = "'43.39"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.27%  '
